# Generate Report for Handback
#
# Re-order the per-file rows on the Overview/zh-cn/de-de sheets so that the
# two files that have been handed back (598b0142... and fe0fa1cb...) move to
# the top, and record their new "Handed back" status plus the new handback
# target-file / handback-file / handback-datetime information.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: remember the external URL behind every existing hyperlink, keyed by
# worksheet name + the *unique* long file name it points at (every hyperlinked
# cell's Address always resolves to a URL that ends with that long file name,
# so we recover the key from the Address itself rather than from
# TextToDisplay -- several cells, e.g. the ".md" column, redundantly display a
# short label for the same link). This lets us restore the correct hyperlink
# target for each file after the rows get shuffled around.
# ---------------------------------------------------------------------------
$urlMap = @{}
foreach ($sheet in $wb.Worksheets) {
    foreach ($hl in $sheet.Hyperlinks) {
        $addr = $hl.Address
        if ([string]::IsNullOrEmpty($addr)) { continue }
        $fileName = $addr.Substring($addr.LastIndexOf("/") + 1)
        $urlMap[$sheet.Name + "|" + $fileName] = $addr
    }
}

# ---------------------------------------------------------------------------
# Step 2: Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ovRows = @{
    2 = @{ A = "598b0142-c9de-47df-9954-dbc2c685e8fa.md"; B = "Handed back: in sync with en-US"; C = "Handed back: in sync with en-US"; D = "2016-17-19 16:17:15" }
    3 = @{ A = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.md"; B = "Handed back: in sync with en-US"; C = "Handed back: in sync with en-US"; D = "2016-17-19 16:17:15" }
    4 = @{ A = "bdac4364-ee81-483c-90f3-a9f81cd0b3a1.md"; B = "In Translation";                  C = "In Translation";                  D = "2016-16-19 16:16:13" }
    5 = @{ A = "9361e3a4-1d7e-466e-a1c0-e8440669a6f9.md"; B = "Ready for handoff";                C = "Ready for handoff";                D = "2016-17-19 16:17:15" }
}

foreach ($r in 2..5) {
    $row = $ovRows[$r]
    $ov.Range("A$r").Value = $row.A
    $ov.Range("B$r").Value = $row.B
    $ov.Range("C$r").Value = $row.C
    $ov.Range("D$r").Value = $row.D
}

foreach ($hl in $ov.Hyperlinks) {
    $ref = $hl.Range.Address($false, $false)
    $rowNum = [int]($ref -replace '[^0-9]', '')
    $fileName = $ovRows[$rowNum].A
    $hl.TextToDisplay = $fileName
    $hl.Address = $urlMap["Overview|" + $fileName]
    # TextToDisplay alone does not refresh the underlying cell text, so set it too.
    $hl.Range.Value = $fileName
}

# ---------------------------------------------------------------------------
# Step 3: per-language detail sheets (zh-cn, de-de)
# ---------------------------------------------------------------------------
$langRows = @{
    "zh-cn" = @{
        2 = @{ A = "598b0142-c9de-47df-9954-dbc2c685e8fa.md"; B = ".md"; C = "Handed back: in sync with en-US";
               D = "598b0142-c9de-47df-9954-dbc2c685e8fa.bbdd8448842836cc1af59cd36da5caab3c9f436d.zh-cn.xlf"; E = "2016-03-19 16:17:12";
               F = "598b0142-c9de-47df-9954-dbc2c685e8fa.md";
               G = "598b0142-c9de-47df-9954-dbc2c685e8fa.bbdd8448842836cc1af59cd36da5caab3c9f436d.zh-cn.xlf";
               H = "2016-03-19 16:17:31"; I = "Include" }
        3 = @{ A = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.md"; B = ".md"; C = "Handed back: in sync with en-US";
               D = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.469b7f37212a929ad20f6d979756347d543be7c2.zh-cn.xlf"; E = "2016-03-19 16:17:12";
               F = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.md";
               G = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.469b7f37212a929ad20f6d979756347d543be7c2.zh-cn.xlf";
               H = "2016-03-19 16:17:31"; I = "Include" }
        4 = @{ A = "bdac4364-ee81-483c-90f3-a9f81cd0b3a1.md"; B = ".md"; C = "In Translation";
               D = "bdac4364-ee81-483c-90f3-a9f81cd0b3a1.efd6bd1465697e78e5b7fa3670a428577a84e725.zh-cn.xlf"; E = "2016-03-19 16:16:10";
               F = $null; G = $null;
               H = "0001-01-01 00:00:00"; I = "Include" }
        5 = @{ A = "9361e3a4-1d7e-466e-a1c0-e8440669a6f9.md"; B = ".md"; C = "Ready for handoff";
               D = "9361e3a4-1d7e-466e-a1c0-e8440669a6f9.f983e3c2d0c160ea20984efc165fcaec781d62f3.zh-cn.xlf"; E = "2016-03-19 16:17:12";
               F = $null; G = $null;
               H = "0001-01-01 00:00:00"; I = "Include" }
    }
    "de-de" = @{
        2 = @{ A = "598b0142-c9de-47df-9954-dbc2c685e8fa.md"; B = ".md"; C = "Handed back: in sync with en-US";
               D = "598b0142-c9de-47df-9954-dbc2c685e8fa.bbdd8448842836cc1af59cd36da5caab3c9f436d.de-de.xlf"; E = "2016-03-19 16:17:15";
               F = "598b0142-c9de-47df-9954-dbc2c685e8fa.md";
               G = "598b0142-c9de-47df-9954-dbc2c685e8fa.bbdd8448842836cc1af59cd36da5caab3c9f436d.de-de.xlf";
               H = "2016-03-19 16:17:35"; I = "Include" }
        3 = @{ A = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.md"; B = ".md"; C = "Handed back: in sync with en-US";
               D = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.469b7f37212a929ad20f6d979756347d543be7c2.de-de.xlf"; E = "2016-03-19 16:17:15";
               F = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.md";
               G = "fe0fa1cb-2894-42ff-97a9-757491a0f6ec.469b7f37212a929ad20f6d979756347d543be7c2.de-de.xlf";
               H = "2016-03-19 16:17:35"; I = "Include" }
        4 = @{ A = "bdac4364-ee81-483c-90f3-a9f81cd0b3a1.md"; B = ".md"; C = "In Translation";
               D = "bdac4364-ee81-483c-90f3-a9f81cd0b3a1.efd6bd1465697e78e5b7fa3670a428577a84e725.de-de.xlf"; E = "2016-03-19 16:16:13";
               F = $null; G = $null;
               H = "0001-01-01 00:00:00"; I = "Include" }
        5 = @{ A = "9361e3a4-1d7e-466e-a1c0-e8440669a6f9.md"; B = ".md"; C = "Ready for handoff";
               D = "9361e3a4-1d7e-466e-a1c0-e8440669a6f9.f983e3c2d0c160ea20984efc165fcaec781d62f3.de-de.xlf"; E = "2016-03-19 16:17:15";
               F = $null; G = $null;
               H = "0001-01-01 00:00:00"; I = "Include" }
    }
}

foreach ($langName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($langName)
    $rows = $langRows[$langName]

    # write the plain (non-hyperlinked) cell values first
    foreach ($r in 2..5) {
        $row = $rows[$r]
        $ws.Range("C$r").Value = $row.C
        $ws.Range("E$r").Value = $row.E
        $ws.Range("H$r").Value = $row.H
        $ws.Range("I$r").Value = $row.I
    }

    # fix up the existing hyperlinked cells (A, B, D) for every row
    foreach ($hl in $ws.Hyperlinks) {
        $ref = $hl.Range.Address($false, $false)
        $col = ($ref -replace '[0-9]', '')
        $rowNum = [int]($ref -replace '[^0-9]', '')
        $row = $rows[$rowNum]

        if ($col -eq "A") {
            $newText = $row.A
            $urlKey = $row.A
        } elseif ($col -eq "B") {
            $newText = $row.B
            $urlKey = $row.A
        } else {
            $newText = $row.D
            $urlKey = $row.D
        }

        $hl.TextToDisplay = $newText
        $hl.Address = $urlMap[$langName + "|" + $urlKey]
        # TextToDisplay alone does not refresh the underlying cell text, so set it too.
        $hl.Range.Value = $newText
    }

    # add the brand-new hyperlinked cells (F, G) for the two handed-back rows
    foreach ($r in 2, 3) {
        $row = $rows[$r]

        $fCell = $ws.Range("F$r")
        $fCell.Value = $row.F
        $fLink = $ws.Hyperlinks.Add($fCell, $urlMap[$langName + "|" + $row.F])
        $fLink.TextToDisplay = $row.F

        $gCell = $ws.Range("G$r")
        $gCell.Value = $row.G
        $gLink = $ws.Hyperlinks.Add($gCell, $urlMap[$langName + "|" + $row.D])
        $gLink.TextToDisplay = $row.G
    }
}
